$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Update the employee id everywhere it appears (shared across both sheets).
# Doing this first means any rows duplicated afterwards already carry the
# correct, new id.
# ---------------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("emp_c2dcy26q", "emp_yde33znx", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
}

# ---------------------------------------------------------------------------
# Sheet 1: "Weekly Timesheet"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Duplicate the "2026-01-02 / Campbell / OT" row so the Tubergen entry can be
# split across a Regular row and an OT row, matching the new full-month
# simulation coverage. Copying (rather than re-typing) keeps the date cell a
# text value instead of Excel re-interpreting it as a serial date.
$ws1.Rows.Item(5).Copy()
$ws1.Rows.Item(6).Insert()

# Row 2 (Hall / Regular) - rate now simulated
$ws1.Cells.Item(2, 5).Value = 92
$ws1.Cells.Item(2, 6).Value = 920

# Row 3 (McGill / Regular) - rate now simulated
$ws1.Cells.Item(3, 5).Value = 92
$ws1.Cells.Item(3, 6).Value = 920

# Row 4 (2026-01-01): Layne/Holiday -> PTO/PTO, hours 20 -> 16
$ws1.Cells.Item(4, 2).Value = "PTO"
$ws1.Cells.Item(4, 3).Value = 16
$ws1.Cells.Item(4, 4).Value = "PTO"
$ws1.Cells.Item(4, 5).Value = 92
$ws1.Cells.Item(4, 6).Value = 1472

# Row 5 (2026-01-02): Campbell/OT -> Tubergen/Regular, hours 20 -> 4
$ws1.Cells.Item(5, 2).Value = "Tubergen"
$ws1.Cells.Item(5, 3).Value = 4
$ws1.Cells.Item(5, 4).Value = "Regular"
$ws1.Cells.Item(5, 5).Value = 92
$ws1.Cells.Item(5, 6).Value = 368

# Row 6 (2026-01-02): Campbell/OT -> Tubergen/OT, hours 20 -> 4
$ws1.Cells.Item(6, 2).Value = "Tubergen"
$ws1.Cells.Item(6, 3).Value = 4
$ws1.Cells.Item(6, 5).Value = 92
$ws1.Cells.Item(6, 6).Value = 552

# Row 8 (was row 7 before the insert): SUBTOTAL
$ws1.Cells.Item(8, 3).Value = 44
$ws1.Cells.Item(8, 4).Value = "Reg: 40 / OT: 4"
$ws1.Cells.Item(8, 6).Value = 4232

# Row 11 (was row 10): HOURLY SUBTOTAL
$ws1.Cells.Item(11, 6).Value = 4232

# Row 13 (was row 12): GRAND TOTAL
$ws1.Cells.Item(13, 6).Value = 4232

# ---------------------------------------------------------------------------
# Sheet 2: "Jason Schema"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Duplicate the Campbell/OT row the same way as sheet 1. Inserting the copy
# ABOVE row 5 (rather than below) means the original row 5 - including its
# trailing blank Notes cell - simply slides down to row 6 intact, while the
# freshly pasted row lands on row 5.
$ws2.Rows.Item(5).Copy()
$ws2.Rows.Item(5).Insert()

# Row 2 (Hall / Regular) - rate now simulated
$ws2.Cells.Item(2, 6).Value = 92
$ws2.Cells.Item(2, 7).Value = 920

# Row 3 (McGill / Regular) - rate now simulated
$ws2.Cells.Item(3, 6).Value = 92
$ws2.Cells.Item(3, 7).Value = 920

# Row 4 (2026-01-01): Layne/Holiday -> PTO/PTO, hours 20 -> 16
$ws2.Cells.Item(4, 4).Value = "PTO"
$ws2.Cells.Item(4, 5).Value = 16
$ws2.Cells.Item(4, 6).Value = 92
$ws2.Cells.Item(4, 7).Value = 1472
$ws2.Cells.Item(4, 8).Value = "PTO"
$ws2.Cells.Item(4, 9).Value = "PTO"

# Row 5 (2026-01-02): Campbell/OT -> Tubergen/Regular, hours 20 -> 4
$ws2.Cells.Item(5, 4).Value = "Tubergen"
$ws2.Cells.Item(5, 5).Value = 4
$ws2.Cells.Item(5, 6).Value = 92
$ws2.Cells.Item(5, 7).Value = 368
$ws2.Cells.Item(5, 8).Value = "Regular"

# Row 6 (2026-01-02): Campbell/OT -> Tubergen/OT, hours 20 -> 4
$ws2.Cells.Item(6, 4).Value = "Tubergen"
$ws2.Cells.Item(6, 5).Value = 4
$ws2.Cells.Item(6, 6).Value = 92
$ws2.Cells.Item(6, 7).Value = 552
